$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above row 856 ("Camote", "1a
# (guarda)"), pushing the existing rows 856..958 down to 857..959.
$ws.Rows.Item(856).Insert()

$ws.Cells.Item(856,1).Value = 3
$ws.Cells.Item(856,2).Value = "Femacal de La Calera"
$ws.Cells.Item(856,3).Value = "Coquimbo"
$ws.Cells.Item(856,4).Value = 45142
$ws.Cells.Item(856,5).Value = 5
$ws.Cells.Item(856,6).Value = 100112045
$ws.Cells.Item(856,7).Value = "Zapallo"
$ws.Cells.Item(856,8).Value = "Camote"
$ws.Cells.Item(856,9).Value = "1a (guarda)"
$ws.Cells.Item(856,10).Value = 120
$ws.Cells.Item(856,11).Value = 500
$ws.Cells.Item(856,12).Value = 500
$ws.Cells.Item(856,13).Value = 500
$ws.Cells.Item(856,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(856,15).Value = "Provincia de Talca"
$ws.Cells.Item(856,16).Value = 500
$ws.Cells.Item(856,17).Value = 1
$ws.Cells.Item(856,18).Value = "Hortaliza"

# Make sure the D column keeps the date number format used by the rest
# of the column (style index 2 in the original workbook).
$ws.Cells.Item(856,4).NumberFormat = $ws.Cells.Item(857,4).NumberFormat
